$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point value in A22 (timestamp rounding fix)
$ws.Range("A22").Value = 45876.83356961806

# Append new row 23 with the latest sensor reading
$ws.Range("A23").Value = 45876.87520870067
$ws.Range("B23").Value = 2025
$ws.Range("C23").Value = 28
$ws.Range("D23").Value = 14.67
$ws.Range("E23").Value = 90.31999999999999
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = "-"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = "21:00:18"

# Match the date/time number format used by the rest of column A
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
